# Adding new test data
# Appends the letter "e" to the end of each test-data string in the
# "reg" worksheet's Username column (C2:C11).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("reg")

$cells = @("C2", "C3", "C4", "C5", "C6", "C7", "C8", "C9", "C10", "C11")

foreach ($addr in $cells) {
    $range = $ws.Range($addr)
    $current = $range.Value()
    $range.Value = [string]$current + "e"
}
